# Updated capital structure database
# Refresh the Falkland Islands Oil/Gas (Production and Exploration)
# capital-structure metrics for both companies (rows 2 and 3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("K2").Value = -0.417
$ws.Range("O2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("U2").Value = 0.5600000000000001
$ws.Range("V2").Value = 0.08641975308641976
$ws.Range("W2").Value = -0.01404040404040404
$ws.Range("X2").Value = 0.09831860657952607
$ws.Range("Y2").Value = -0.1123590106199301
$ws.Range("AA2").Value = -0.0140406695255222
$ws.Range("AB2").Value = 0.09831860657952607
$ws.Range("AC2").Value = -0.1123592761050483
$ws.Range("AG2").Value = -0.5600000000000001
$ws.Range("AJ2").Value = -0.0945945945945946
$ws.Range("AK2").Value = -0.01948503827418233
$ws.Range("AM2").Value = -0.003
$ws.Range("AQ2").Value = 135.3333333333333

# --- Row 3 ---
$ws.Range("K3").Value = -0.417
$ws.Range("O3").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("U3").Value = 0.5600000000000001
$ws.Range("V3").Value = 0.08641975308641976
$ws.Range("W3").Value = -0.01404040404040404
$ws.Range("X3").Value = 0.09831860657952607
$ws.Range("Y3").Value = -0.1123590106199301
$ws.Range("AA3").Value = -0.0140406695255222
$ws.Range("AB3").Value = 0.09831860657952607
$ws.Range("AC3").Value = -0.1123592761050483
$ws.Range("AG3").Value = -0.5600000000000001
$ws.Range("AJ3").Value = -0.0945945945945946
$ws.Range("AK3").Value = -0.01948503827418233
$ws.Range("AM3").Value = -0.003
$ws.Range("AQ3").Value = 135.3333333333333
